# Added new TC TC_Name_2
# - Test_Cases sheet: Test_Case_04 renamed to RegisterUserForTheSeleniumTraining2, Run_Mode -> Y
# - Test_Data sheet: new "RegisterUserForTheSeleniumTraining2" section appended (rows 19-21),
#   mirroring the first ("RegisterUserForTheSeleniumTraining") section's layout/content.
# - Active sheet switched from Test_Cases to Test_Data, with new selections on each sheet.

$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("Test_Cases")
$wsData  = $wb.Worksheets.Item("Test_Data")

# ---------------------------------------------------------------------------
# 1) Test_Cases: row 5 becomes the new test case, enabled for run.
# ---------------------------------------------------------------------------
$wsCases.Range("B5").Value = "RegisterUserForTheSeleniumTraining2"
$wsCases.Range("C5").Value = "Y"

# ---------------------------------------------------------------------------
# 2) Test_Data: append a 4th data block (header + column titles + Pallav row),
#    cloned from the first block (rows 1-3) so formatting/styles match.
# ---------------------------------------------------------------------------
$wsData.Range("A1").Copy()
$wsData.Range("A19").PasteSpecial(-4122)

$wsData.Range("A2:G2").Copy()
$wsData.Range("A20").PasteSpecial(-4122)

$wsData.Range("A3:G3").Copy()
$wsData.Range("A21").PasteSpecial(-4122)

$wsData.Range("A19").Value = "RegisterUserForTheSeleniumTraining2"

$wsData.Range("A20").Value = "Sr.No"
$wsData.Range("B20").Value = "First_Name"
$wsData.Range("C20").Value = "Email_ID"
$wsData.Range("D20").Value = "Phone_No"
$wsData.Range("E20").Value = "Message"
$wsData.Range("F20").Value = "Password"
$wsData.Range("G20").Value = "Page_Title"

$wsData.Range("A21").Value = "'1"
$wsData.Range("B21").Value = "Pallav"
$wsData.Range("C21").Value = "Pallav@gmail.com"
$wsData.Range("D21").Value = "'353535353"
$wsData.Range("E21").Value = "Message_Pallav"
$wsData.Range("F21").Value = "passd1235"
$wsData.Range("G21").Value = "My Store"

$wsData.Hyperlinks.Add($wsData.Range("C21"), "mailto:Pallav@gmail.com")

# Re-stamp C21's number format/border/etc. (Hyperlinks.Add recolors the cell
# with the hyperlink style) so it matches the other mailto cells again.
$wsData.Range("C3").Copy()
$wsData.Range("C21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Selections / active sheet, matching the saved UI state in the workbook.
# ---------------------------------------------------------------------------
[void]$wsCases.Range("C15").Select()
[void]$wsData.Range("A20").Select()
[void]$wsData.Activate()
